# Lab 4 - Entrega Final
# The "Datos Lab4" sheet holds two result blocks (ARRAYLIST @ rows 2-11,
# LINKED_LIST @ rows 15-24), each with 10 sample sizes (1000..512000).
# The final run only kept the first 4 sample sizes (1000, 2000, 4000, 8000)
# per block, and the measured times were refreshed with the latest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab4")

function Set-Block($range, $values) {
    $rows = $values.Count
    $cols = $values[0].Count
    $arr = New-Object 'object[,]' $rows, $cols
    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $arr[$r, $c] = $values[$r][$c]
        }
    }
    $ws.Range($range).Value = $arr
}

# --- Drop the rows for the sample sizes that are no longer reported ---
# (bottom-up so row numbers don't shift under us)
$ws.Range("A19:A24").EntireRow.Delete() | Out-Null
$ws.Range("A6:A11").EntireRow.Delete() | Out-Null

# --- Refresh the measured values for the ARRAYLIST block (rows 2:5) ---
Set-Block "B2:D5" @(
    @(5406.25, 5375, 5531.25),
    @(29109.375, 38968.75, 98921.875),
    @(139640.625, 140953.125, 140109.375),
    @(662125, 661968.75, 662437.5)
)

# --- Refresh the measured values for the LINKED_LIST block (rows 9:12) ---
Set-Block "B9:D12" @(
    @(5359.375, 5359.375, 5531.25),
    @(29578.125, 29765.625, 28984.375),
    @(140343.375, 142187.5, 141359.375),
    @(662022.35, 663198.25, 660112.5)
)

# --- Match the saved selection/active cell ---
$ws.Range("B14").Select()
